# Applies the commit "New papers on adoption of US, writing, etc" to
# References/VisualizationSurvey.docx:
#  - Inserts a new "Tran2017 - Beyond Ultrasound Guidance for Regional
#    Anesthesiology:" section (heading + 4 bullets) at the very end of the
#    document (after the Brudfors2015 bullets, before the sectPr).
#  - Moves the "_GoBack" bookmark from the end of the Nagpal2015 section
#    to the end of the new, final paragraph.
#  - A fresh bulleted-list definition (abstractNum) is minted for the new
#    bullets, as Word does whenever a "new" list is created; the existing
#    abstractNum/num entries are renumbered to make room for it, exactly
#    as Word would on save.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Pull the whole package (document.xml + numbering.xml, etc.) as one
#    OOXML string so we can perform the structural edits with ordinary
#    string/regex operations, then push the result back in one shot.
# ---------------------------------------------------------------------
$full = $d.Content.WordOpenXML

# ---------------------------------------------------------------------
# 2. document.xml: relocate the _GoBack bookmark and append the new
#    "Tran2017" section.
# ---------------------------------------------------------------------
$bookmark = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
if ($full.IndexOf($bookmark) -lt 0) {
    throw "could not find _GoBack bookmark"
}
$full = $full.Replace($bookmark, '')

$tran = ''
$tran += '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr>'
$tran += '<w:r><w:rPr><w:b/></w:rPr><w:t>Tran2017 – Beyond Ultrasound Guidance for Regional Anesthesiology:</w:t></w:r>'
$tran += '</w:p>'
$tran += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr>'
$tran += '<w:r><w:t>Concise review of the US in regional anesthesiology including uses, benefits, limitations, and improvements</w:t></w:r>'
$tran += '</w:p>'
$tran += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr>'
$tran += '<w:r><w:lastRenderedPageBreak/><w:t>UGRA generally demonstrates benefits in procedure time, needle insertion attempts, and local anesthetic system toxicity (LAST) occurrence rates versus palpation or otherwise not using ultrasound. Results are most pronounced for inexperienced operators or difficult spinal anatomy.</w:t></w:r>'
$tran += '</w:p>'
$tran += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr>'
$tran += '<w:r><w:t>Limitations of UGRA include difficulty in visualizing complex skeletal and deep nervous structures</w:t></w:r>'
$tran += '<w:r><w:t>.</w:t></w:r>'
$tran += '</w:p>'
$tran += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr>'
$tran += '<w:r><w:t>Suggested improvements to the current state of the art of US pertain largely to operators</w:t></w:r>'
$tran += '<w:r><w:t xml:space="preserve"> and authors</w:t></w:r>'
$tran += '<w:r><w:t xml:space="preserve">, and include using other tools such as fluoroscopy when appropriate, </w:t></w:r>'
$tran += '<w:r><w:t>and a focusing of publication efforts as opposed to the current practice of publishing every variation of each method as a “Brief Technical Report”.</w:t></w:r>'
$tran += $bookmark
$tran += '</w:p>'

$sectPrIdx = $full.IndexOf('<w:sectPr')
if ($sectPrIdx -lt 0) {
    throw "could not find sectPr"
}
$full = $full.Substring(0, $sectPrIdx) + $tran + $full.Substring($sectPrIdx)

# ---------------------------------------------------------------------
# 3. numbering.xml: mint a new abstractNum (bullet list) for numId=9,
#    shifting the existing abstractNumId=2..7 up to 3..8 to make room,
#    mirroring what Word does when a brand-new list is first used.
# ---------------------------------------------------------------------
foreach ($pair in @(@(7,8), @(6,7), @(5,6), @(4,5), @(3,4), @(2,3))) {
    $oldId = $pair[0]
    $newId = $pair[1]
    $full = $full.Replace('w:abstractNumId="' + $oldId + '"', 'w:abstractNumId="' + $newId + '"')
    $full = $full.Replace('<w:abstractNumId w:val="' + $oldId + '"/>', '<w:abstractNumId w:val="' + $newId + '"/>')
}

# Clone the (untouched) abstractNumId="1" block as the template for the
# new list definition, then re-key it to a fresh nsid/tmpl pair.
$templatePattern = '<w:abstractNum w:abstractNumId="1"[^>]*>.*?</w:abstractNum>'
$templateMatch = [regex]::Match($full, $templatePattern, [System.Text.RegularExpressions.RegexOptions]::Singleline)
if (-not $templateMatch.Success) {
    throw "could not find abstractNum template"
}
$newAbstractNum = $templateMatch.Value
$newAbstractNum = $newAbstractNum.Replace('w:abstractNumId="1"', 'w:abstractNumId="2"')
$newAbstractNum = $newAbstractNum.Replace('<w:nsid w:val="1F355D8E"/>', '<w:nsid w:val="1F570593"/>')
$newAbstractNum = $newAbstractNum.Replace('<w:tmpl w:val="905214D4"/>', '<w:tmpl w:val="F176F678"/>')

$insertBefore = '<w:abstractNum w:abstractNumId="3" w15:restartNumberingAfterBreak="0"><w:nsid w:val="54690072"/>'
$insertIdx = $full.IndexOf($insertBefore)
if ($insertIdx -lt 0) {
    throw "could not find abstractNum insertion point"
}
$full = $full.Substring(0, $insertIdx) + $newAbstractNum + $full.Substring($insertIdx)

# New <w:num> entry mapping numId=9 to the freshly minted abstractNum.
$newNum = '<w:num w:numId="9"><w:abstractNumId w:val="2"/></w:num>'
$numberingEnd = '</w:numbering>'
$endIdx = $full.LastIndexOf($numberingEnd)
if ($endIdx -lt 0) {
    throw "could not find numbering end"
}
$full = $full.Substring(0, $endIdx) + $newNum + $full.Substring($endIdx)

# ---------------------------------------------------------------------
# 4. Push the rewritten package back into the document.
# ---------------------------------------------------------------------
$d.Content.InsertXML($full)
